$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(5, 6).Value = 80
$ws.Cells.Item(6, 6).Value = 861
$ws.Cells.Item(7, 6).Value = 442
$ws.Cells.Item(8, 6).Value = 4757
$ws.Cells.Item(9, 6).Value = 4757
$ws.Cells.Item(12, 6).Value = 163
$ws.Cells.Item(14, 6).Value = 201
$ws.Cells.Item(16, 6).Value = 7600
$ws.Cells.Item(17, 6).Value = 252
$ws.Cells.Item(18, 6).Value = 130
$ws.Cells.Item(19, 6).Value = 296
$ws.Cells.Item(21, 6).Value = 534
$ws.Cells.Item(22, 6).Value = 1400
$ws.Cells.Item(25, 6).Value = 2258
$ws.Cells.Item(27, 6).Value = 2093
$ws.Cells.Item(28, 6).Value = 6197
$ws.Cells.Item(34, 6).Value = 6492
$ws.Cells.Item(39, 6).Value = 22
$ws.Cells.Item(41, 6).Value = 2462
$ws.Cells.Item(46, 6).Value = 452
$ws.Cells.Item(47, 6).Value = 2153
$ws.Cells.Item(49, 6).Value = 1084

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 234

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 1453

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 1453
$ws.Cells.Item(5, 6).Value = 234
$ws.Cells.Item(6, 6).Value = 80
$ws.Cells.Item(8, 6).Value = 442
$ws.Cells.Item(9, 6).Value = 4757
$ws.Cells.Item(10, 6).Value = 4757
$ws.Cells.Item(13, 6).Value = 163
$ws.Cells.Item(15, 6).Value = 201
$ws.Cells.Item(17, 6).Value = 7600
$ws.Cells.Item(18, 6).Value = 252
$ws.Cells.Item(19, 6).Value = 130
$ws.Cells.Item(20, 6).Value = 534
$ws.Cells.Item(21, 6).Value = 1400
$ws.Cells.Item(24, 6).Value = 2258
$ws.Cells.Item(26, 6).Value = 2093
$ws.Cells.Item(29, 6).Value = 6197
$ws.Cells.Item(35, 6).Value = 448
$ws.Cells.Item(36, 6).Value = 6492
$ws.Cells.Item(40, 6).Value = 22
$ws.Cells.Item(42, 6).Value = 2462
$ws.Cells.Item(46, 6).Value = 452
$ws.Cells.Item(48, 6).Value = 2153
